$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix up EXECUTE / PARALLEL flags on existing rows ---
$ws.Range("H3").Value() = "No"
$ws.Range("B5").Value() = "YES"
$ws.Range("H6").Value() = "YES"
$ws.Range("H7").Value() = "YES"
$ws.Range("B8").Value() = "YES"
$ws.Range("H8").Value() = "YES"
$ws.Range("B9").Value() = "YES"
$ws.Range("H9").Value() = "YES"

# --- Insert 4 new test case rows before the API test cases ---
$ws.Rows("10:13").Insert()

$ws.Range("A10").Value() = "TC_009_GO_TO_HELP_PAGE"
$ws.Range("B10").Value() = "YES"
$ws.Range("C10").Value() = "Web"
$ws.Range("D10").Value() = "TestHomePage"
$ws.Range("E10").Value() = "go_to_help"
$ws.Range("F10").Value() = "DATAFILE.xlsx"
$ws.Range("G10").Value() = "DATA_FILE"
$ws.Range("H10").Value() = "YES"

$ws.Range("A11").Value() = "TC_010_Go_TO_CONTACT_US_PAGE"
$ws.Range("B11").Value() = "No"
$ws.Range("C11").Value() = "Web"
$ws.Range("D11").Value() = "TestHomePage"
$ws.Range("E11").Value() = "go_to_contact_us_page"
$ws.Range("F11").Value() = "DATAFILE.xlsx"
$ws.Range("G11").Value() = "DATA_FILE"
$ws.Range("H11").Value() = "YES"

$ws.Range("A12").Value() = "TC_011_CHANGE_COUNTRY"
$ws.Range("B12").Value() = "YES"
$ws.Range("C12").Value() = "Web"
$ws.Range("D12").Value() = "TestHomePage"
$ws.Range("E12").Value() = "test_change_country"
$ws.Range("F12").Value() = "DATAFILE.xlsx"
$ws.Range("G12").Value() = "DATA_FILE"
$ws.Range("H12").Value() = "YES"

$ws.Range("A13").Value() = "TC_012_ADD_NEW_ADDRESS_PAGE"
$ws.Range("B13").Value() = "No"
$ws.Range("C13").Value() = "Web"
$ws.Range("D13").Value() = "TestHomePage"
$ws.Range("E13").Value() = "add_new_address"
$ws.Range("F13").Value() = "DATAFILE.xlsx"
$ws.Range("G13").Value() = "DATA_FILE"
$ws.Range("H13").Value() = "YES"

# --- Update selection to match the authored state ---
$ws.Range("B13").Select()
